# Auto-generated edit script implementing the diff for sheet 'ランサーズ'
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Remove the now-stale hyperlinks for the whole sheet; they will be
#        re-created below only for the rows that survive (F2:F9). Deleting
#        individual Hyperlinks.Item(i) / Range.Hyperlinks.Delete() does not
#        reliably scope to one cell in this host, so clear-and-rebuild is used.
$ws.Cells.Hyperlinks.Delete()

# --- 2. Delete the rows that no longer exist in the refreshed scrape (10-21).
$ws.Range("A10:H21").EntireRow.Delete()

# --- 3. Overwrite rows 2-9 with the refreshed listing data.

# Row 2
$ws.Range("A2").Value2 = '2025-12-10 06:30:13'
$ws.Range("B2").Value2 = '産業機械向けAI異常検知・状態推定システムの開発・導入支援エンジニア募集(AI/エッジ・組み込み)'
$ws.Range("C2").Value2 = 'システム開発'
$ws.Range("D2").Value2 = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E2").Value2 = '期限情報なし'
$ws.Range("F2").Value2 = 'https://www.lancers.jp/work/detail/5450864'
$ws.Range("G2").Value2 = 383
$ws.Range("H2").Value2 = '🔥AI,Ai ◆開発'

# Row 3
$ws.Range("A3").Value2 = '2025-12-10 06:30:13'
$ws.Range("B3").Value2 = '【自動化】Webサービス更新ツール開発(200アカウント管理)'
$ws.Range("C3").Value2 = 'システム開発'
$ws.Range("D3").Value2 = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E3").Value2 = '期限情報なし'
$ws.Range("F3").Value2 = 'https://www.lancers.jp/work/detail/5448409'
$ws.Range("G3").Value2 = 230
$ws.Range("H3").Value2 = '◆ツール,開発 ◇管理'

# Row 4
$ws.Range("A4").Value2 = '2025-12-10 06:30:13'
$ws.Range("B4").Value2 = '【急募】某新聞社のプロトタイプシステム用チャットボット開発'
$ws.Range("C4").Value2 = 'システム開発'
$ws.Range("D4").Value2 = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E4").Value2 = '期限情報なし'
$ws.Range("F4").Value2 = 'https://www.lancers.jp/work/detail/5450641'
$ws.Range("G4").Value2 = 83
$ws.Range("H4").Value2 = '◆開発'

# Row 5
$ws.Range("A5").Value2 = '2025-12-10 06:30:13'
$ws.Range("B5").Value2 = '在宅専業OK│フルスタックエンジニア/開発×データ処理に挑戦!EC運営を支える仕事!'
$ws.Range("C5").Value2 = 'システム開発'
$ws.Range("D5").Value2 = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E5").Value2 = '期限情報なし'
$ws.Range("F5").Value2 = 'https://www.lancers.jp/work/detail/5450846'
$ws.Range("G5").Value2 = 75
$ws.Range("H5").Value2 = '◆開発'

# Row 6
$ws.Range("A6").Value2 = '2025-12-10 06:30:13'
$ws.Range("B6").Value2 = '【フルスタックエンジニア募集】新規Webサービス開発'
$ws.Range("C6").Value2 = 'システム開発'
$ws.Range("D6").Value2 = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E6").Value2 = '期限情報なし'
$ws.Range("F6").Value2 = 'https://www.lancers.jp/work/detail/5450548'
$ws.Range("G6").Value2 = 75
$ws.Range("H6").Value2 = '◆開発'

# Row 7
$ws.Range("A7").Value2 = '2025-12-10 06:30:13'
$ws.Range("B7").Value2 = '【若手歓迎×リモートOK】SRE/インフラエンジニア(Google Cloud/長期・金融系案件)'
$ws.Range("C7").Value2 = 'システム開発'
$ws.Range("D7").Value2 = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E7").Value2 = '期限情報なし'
$ws.Range("F7").Value2 = 'https://www.lancers.jp/work/detail/5445466'
$ws.Range("G7").Value2 = 25
$ws.Range("H7").ClearContents()

# Row 8
$ws.Range("A8").Value2 = '2025-12-10 06:30:13'
$ws.Range("B8").Value2 = '注目 限定公開 PR 限定公開の仕事'
$ws.Range("C8").Value2 = 'システム開発'
$ws.Range("D8").Value2 = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E8").Value2 = '期限情報なし'
$ws.Range("F8").Value2 = 'https://www.lancers.jp/work/detail/5450323'
$ws.Range("G8").Value2 = 13
$ws.Range("H8").ClearContents()

# Row 9
$ws.Range("A9").Value2 = '2025-12-10 06:30:13'
$ws.Range("B9").Value2 = '【急募】当社HPのバグ修正をお願いしたいです'
$ws.Range("C9").Value2 = 'システム開発'
$ws.Range("D9").Value2 = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E9").Value2 = '期限情報なし'
$ws.Range("F9").Value2 = 'https://www.lancers.jp/work/detail/5450784'
$ws.Range("G9").Value2 = 10
$ws.Range("H9").ClearContents()

# --- 4. Re-create the hyperlinks on column F for the surviving rows, in order
#        (so the relationship ids come out as rId1..rId8, same as the diff).
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5450864', "", "", 'https://www.lancers.jp/work/detail/5450864') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5448409', "", "", 'https://www.lancers.jp/work/detail/5448409') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5450641', "", "", 'https://www.lancers.jp/work/detail/5450641') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5450846', "", "", 'https://www.lancers.jp/work/detail/5450846') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5450548', "", "", 'https://www.lancers.jp/work/detail/5450548') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5445466', "", "", 'https://www.lancers.jp/work/detail/5445466') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5450323', "", "", 'https://www.lancers.jp/work/detail/5450323') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5450784', "", "", 'https://www.lancers.jp/work/detail/5450784') | Out-Null

# --- 5. Column width tweaks: D 32 -> 30, H 14 -> 13 (COM ColumnWidth reports
#        0.83 narrower than the saved OOXML 'width' attribute on this host,
#        so subtract that fixed offset to land exactly on the target widths).
$ws.Columns("D").ColumnWidth = 29.17
$ws.Columns("H").ColumnWidth = 12.17

# --- 6. Dimension (A1:H9) is recomputed automatically by the host on save
#        once the trailing rows/cells are gone.
